$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '60.850.49'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.592.38'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.08%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '522.98'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +2.95%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.08'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +2.42%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.70'
$ws.Range('D9').Style = "Normal"
$ws.Range('E10').Value = '  +1.10%  '
$ws.Range('E11').Value = '  -0.46%  '
$ws.Range('E12').Value = '  +1.31%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '3.046.86'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +0.04%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '60.867.46'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '21.64'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -0.01%  '
$ws.Range('E16').Value = '  +0.05%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.593.53'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.22%  '
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '353.18'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.06%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '10.57'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.22%  '
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '60.79'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.27%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '0.427'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('E25').Value = '  -0.67%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.708.26'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.25%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.00'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.23%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.0₃0844'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.30%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.37'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '6.32'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +10.27%  '
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('E33').Value = '  +2.61%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '148.03'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -3.66%  '
$ws.Range('E35').Value = '  +4.10%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.935'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +9.56%  '
$ws.Range('E37').Value = '  +0.65%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.860'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.50'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  +1.60%  '
$ws.Range('E40').Value = '  +1.20%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '36.45'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.13%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '288.45'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.56%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.102'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +1.86%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.620'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.37%  '
$ws.Range('E45').Value = '  +0.40%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.997'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -0.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '19.58'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.0238'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +1.96%  '
$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '4.88'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '10.32'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.12%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.10'
$ws.Range('D51').Style = "Normal"
